$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67, shifting existing rows 67-190 down to 68-191.
$ws.Rows(67).Insert()

# Populate the newly inserted row 67 with the new weekly record.
$ws.Range("A67").Value = 11
$ws.Range("B67").Value = "Vega Monumental Concepción"
$ws.Range("C67").Value = "Bíobío"
$ws.Range("D67").Value = 44775
$ws.Range("D67").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E67").Value = 8
$ws.Range("F67").Value = 100112003
$ws.Range("G67").Value = "Ajo"
$ws.Range("H67").Value = "Chino"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 300
$ws.Range("K67").Value = 23000
$ws.Range("L67").Value = 24000
$ws.Range("M67").Value = 23500
$ws.Range("N67").Value = "`$/caja 10 kilos"
$ws.Range("O67").Value = "China"
$ws.Range("P67").Value = 2350
$ws.Range("Q67").Value = 10
$ws.Range("R67").Value = "Hortaliza"
